{"js": "// The minutes text originally read \"... reelect the new scrum master ...\" and\n// \"... responsibility to be the new scrum master. He ...\". The final edit\n// drops the now-redundant word \"new\" in both spots (Abdullah was already\n// elected, so he isn't a \"new\" scrum master anymore).\nconst body = context.document.body;\nconst results = body.search(\"new scrum master\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"scrum master\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The minutes text originally read \"... reelect the new scrum master ...\" and\n# \"... responsibility to be the new scrum master. He ...\". The final edit\n# drops the now-redundant word \"new\" in both spots (Abdullah was already\n# elected, so he isn't a \"new\" scrum master anymore).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"new scrum master\"\n$find.Replacement.Text = \"scrum master\"\n$find.Forward = $true\n$find.Wrap = 1          # wdFindContinue - keep searching the whole story\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.MatchSoundsLike = $false\n$find.MatchAllWordForms = $false\n\n# wdReplaceAll = 2 -> replace every occurrence in the range.\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $find.MatchSoundsLike, $find.MatchAllWordForms, $find.Forward, $find.Wrap, $find.Format, $find.Replacement.Text, 2) | Out-Null\n"}
